# menuExcToHtml.xlsx — add a new "Pizze del Mese" (pizza-of-the-month) rotation
# table below the existing data (rows 117-149), reusing the same "section title"
# look already used for Le Classiche / Le Speciali / etc. (rows 23, 45, 94, 99, 109).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Section title "Pizze del Mese" in B117 (merged B117:D118), styled the same
#    way as the other section headers (e.g. B23:D24 "Le Classiche").
# ---------------------------------------------------------------------------
$ws.Range("B117").Value = "Pizze del Mese"

# Merge first, then paste formatting from a single already-styled "section
# header" cell onto the merged block as one unit, so every cell in the new
# block shares the exact same uniform style (no per-edge merge borders).
$ws.Range("B117:D118").Merge()
$ws.Range("B23").Copy()
$ws.Range("B117:D118").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B117").Value = "Pizze del Mese"

# ---------------------------------------------------------------------------
# 2. Weekly rotation list: col A = week number, col B = month name.
#    Fill these in row order first so the new shared-strings land in the same
#    order the month names are first encountered (repeats reuse the string).
# ---------------------------------------------------------------------------
$months = @(
  @(119, 18, "Novembre"),
  @(120, 18, "Dicembre"),
  @(121, 19, "Gennaio"),
  @(122, 19, "Febbraio"),
  @(123, 19, "Marzo"),
  @(124, 19, "Aprile"),
  @(125, 19, "Maggio"),
  @(126, 19, "Giugno"),
  @(127, 19, "Luglio"),
  @(128, 19, "Agosto"),
  @(129, 19, "Settembre"),
  @(130, 19, "Ottobre"),
  @(131, 19, "Novembre"),
  @(132, 19, "Dicembre"),
  @(133, 20, "Gennaio"),
  @(134, 20, "Febbraio"),
  @(135, 20, "Marzo"),
  @(136, 20, "Aprile"),
  @(137, 20, "Maggio"),
  @(138, 20, "Giugno"),
  @(139, 20, "Luglio"),
  @(140, 20, "Agosto"),
  @(141, 20, "Settembre"),
  @(142, 20, "Ottobre"),
  @(143, 20, "Novembre"),
  @(144, 20, "Dicembre"),
  @(145, 21, "Gennaio"),
  @(146, 21, "Febbraio"),
  @(147, 21, "Marzo"),
  @(148, 21, "Aprile"),
  @(149, $null, "Maggio")
)

foreach ($row in $months) {
    $r = $row[0]
    $week = $row[1]
    $month = $row[2]
    if ($week -ne $null) {
        $ws.Cells.Item($r, 1).Value = $week
    }
    $ws.Cells.Item($r, 2).Value = $month
}

# ---------------------------------------------------------------------------
# 3. Two new pizza-of-the-month descriptions in column C, entered in this
#    specific order (C122 before C121) to match the authoring order.
# ---------------------------------------------------------------------------
$ws.Cells.Item(122, 3).Value = "BASE CREMA DI BROCCOLI: Mozzarella - Scamorza - Lardo in concia"
$ws.Cells.Item(121, 3).Value = "BASE MARGHERITA: Cardi di Cervia - Brie - Alici"

# ---------------------------------------------------------------------------
# 4. Update the view: scroll down to the new section and select C123.
# ---------------------------------------------------------------------------
$ws.Range("C123").Select()
$excel.ActiveWindow.ScrollRow = 112
$excel.ActiveWindow.ScrollColumn = 1
